$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changes from 46076 to 46077 for every data row (rows 2-16)
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 46077
}

# Rows 6-16 are re-populated (reordered) with new A (Beteckning), B (Datum) and G (Area (ha)) values
$rowsData = @(
    @{ Row = 6;  A = "A 5792-2024";  B = 45335;              G = 5.6 },
    @{ Row = 7;  A = "A 7333-2025";  B = 45703.35899305555;  G = 0.9 },
    @{ Row = 8;  A = "A 35642-2023"; B = 45147;              G = 1.2 },
    @{ Row = 9;  A = "A 28288-2023"; B = 45099.6349537037;   G = 0.5 },
    @{ Row = 10; A = "A 12651-2022"; B = 44641;              G = 3.2 },
    @{ Row = 11; A = "A 2593-2024";  B = 45313.69204861111;  G = 2.3 },
    @{ Row = 12; A = "A 13651-2023"; B = 45006;              G = 2.2 },
    @{ Row = 13; A = "A 8194-2025";  B = 45708;              G = 1.9 },
    @{ Row = 14; A = "A 50997-2025"; B = 45946;              G = 1.5 },
    @{ Row = 15; A = "A 7814-2026";  B = 46062.61388888889;  G = 1.1 },
    @{ Row = 16; A = "A 7827-2026";  B = 46062.63958333333;  G = 2.1 }
)

foreach ($item in $rowsData) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 7).Value = $item.G
}
